# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price cells that would otherwise be auto-parsed by Excel as numbers are
# forced to Text format first so they stay literal strings, matching the
# original inlineStr cells (e.g. "1.00", "0.151", "23.92").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.852.77'
$ws.Range("E2").Value = '  -6.64%  '
$ws.Range("D3").Value = '3.296.23'
$ws.Range("E3").Value = '  -5.20%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.55'
$ws.Range("E5").Value = '  -4.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.30'
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.297.15'
$ws.Range("E8").Value = '  -5.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.469'
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("E10").Value = '  -4.35%  '
$ws.Range("E11").Value = '  -5.82%  '
$ws.Range("E12").Value = '  -4.21%  '
$ws.Range("D13").Value = '3.859.87'
$ws.Range("E13").Value = '  -5.09%  '
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '3.300.89'
$ws.Range("E15").Value = '  -5.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000165'
$ws.Range("E16").Value = '  -6.67%  '
$ws.Range("D17").Value = '60.132.24'
$ws.Range("E17").Value = '  -6.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '23.92'
$ws.Range("E18").Value = '  -4.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.59'
$ws.Range("E19").Value = '  -1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.18'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.81'
$ws.Range("E21").Value = '  -11.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '348.73'
$ws.Range("E22").Value = '  -9.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.550'
$ws.Range("E23").Value = '  -2.95%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '3.426.47'
$ws.Range("E25").Value = '  -5.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '68.39'
$ws.Range("E26").Value = '  -8.19%  '
$ws.Range("E27").Value = '  -4.08%  '
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.27'
$ws.Range("E29").Value = '  +2.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.43'
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.77'
$ws.Range("E31").Value = '  -2.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.151'
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("E33").Value = '  -6.43%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = '3.323.24'
$ws.Range("E35").Value = '  -5.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.58'
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.27'
$ws.Range("E37").Value = '  +1.32%  '
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.46'
$ws.Range("E39").Value = '  -2.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '156.51'
$ws.Range("E40").Value = '  -3.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0744'
$ws.Range("E41").Value = '  -4.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.29'
$ws.Range("E43").Value = '  -2.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.27'
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.738'
$ws.Range("E45").Value = '  -7.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.15'
$ws.Range("E46").Value = '  +2.46%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.47'
$ws.Range("E47").Value = '  -4.50%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.53'
$ws.Range("E48").Value = '  -5.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.68'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.75'
$ws.Range("E50").Value = '  +6.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.845'
$ws.Range("E51").Value = '  -5.88%  '
